$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update data: product catalog (Plate / Bench for crunches / Gripper / Water bottle)
$data = @(
    @("Plate", 85.5, "30kg - Steel stainless"),
    @("Plate", 25, "10kg - stainless steel "),
    @("Bench for crunches", 85.12, "Leather bench made specifically for crunches and similar activities"),
    @("Gripper", 25, "5kg resistance"),
    @("Gripper", 35.5, "10kg resistance"),
    @("Water bottle", 12.95, "Plastic water bottle with clapper cap")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Column widths (target stored widths are 22.42578125 / 57.85546875; the
# COM ColumnWidth setter here quantizes to a 1/6-character grid, so feed it
# the input that lands closest to those stored values: 22.5 and 57.8333...)
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(3).ColumnWidth = 57

# Selection
$ws.Range("C16").Select()
